# Apply scheduled-runner market price/profit updates to the Sheets workbook.
# For each leve row listed below, refresh price columns (H-N) with the
# latest computed market values. Cells that no longer apply are cleared;
# new cells are added where a value now exists.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2: Quicksilver
$ws.Range("H2").Value = 557.6
$ws.Range("I2").Value = 129.66667
$ws.Range("J2").Value = 741
$ws.Range("K2").Value = 129.66667
$ws.Range("L2").Value = 741
$ws.Range("M2").Value = -16.66667000000001
$ws.Range("N2").Value = -967
# Row 9: Distilled Water
$ws.Range("H9").Value = 110.42857
$ws.Range("I9").Value = 87.75
$ws.Range("J9").Value = 140.66667
$ws.Range("K9").Value = 87.75
$ws.Range("L9").Value = 140.66667
$ws.Range("M9").Value = 81.25
$ws.Range("N9").Value = -478.66667
# Row 12: Beeswax
$ws.Range("H12").Value = 217.4
$ws.Range("I12").Value = 197.71428
$ws.Range("K12").Value = 197.71428
$ws.Range("M12").Value = -27.71428
# Row 32: Crab Oil
$ws.Range("H32").Value = 26665.666
$ws.Range("I32").Value = 26665.666
$ws.Range("K32").Value = 26665.666
$ws.Range("M32").Value = -26339.666
# Row 38: Hi-Potion of Strength
$ws.Range("H38").Value = 814.26666
$ws.Range("I38").Value = 449.8
$ws.Range("J38").Value = 996.5
$ws.Range("K38").Value = 1349.4
$ws.Range("L38").Value = 2989.5
$ws.Range("M38").Value = -977.4000000000001
$ws.Range("N38").Value = -3733.5
# Row 40: Horn Glue
$ws.Range("H40").Value = 2498
$ws.Range("I40").Value = 1997
$ws.Range("K40").Value = 1997
$ws.Range("M40").Value = -1822
# Row 51: Shark Oil
$ws.Range("H51").Value = 259748.25
$ws.Range("I51").Value = 12998
$ws.Range("J51").Value = 999999
$ws.Range("K51").Value = 12998
$ws.Range("L51").Value = 999999
$ws.Range("M51").Value = -12514
$ws.Range("N51").Value = -1000967
# Row 58: Mega-Potion of Vitality
$ws.Range("H58").Value = 1617
$ws.Range("I58").Value = 449.23077
$ws.Range("J58").Value = 2629.0667
$ws.Range("K58").Value = 1347.69231
$ws.Range("L58").Value = 7887.2001
$ws.Range("M58").Value = -1197.69231
$ws.Range("N58").Value = -8187.2001
# Row 62: Enchanted Mythrite Ink
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
# Row 65: Enchanted Mythrite Ink
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
# Row 80: Hallowed Water
$ws.Range("H80").Value = 623.875
$ws.Range("I80").Value = 573
$ws.Range("K80").Value = 1719
$ws.Range("M80").Value = -721
# Row 83: Hallowed Water
$ws.Range("H83").Value = 623.875
$ws.Range("I83").Value = 573
$ws.Range("K83").Value = 5157
$ws.Range("M83").Value = -165
# Row 87: Noble Gold
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
# Row 90: Noble Gold
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
# Row 137: Magnesia Whetstone
$ws.Range("H137").Value = 1999
$ws.Range("I137").Value = 1999
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 5997
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -3447
$ws.Range("N137").ClearContents()
# Row 141: Grade 1 Gemdraught of Mind
$ws.Range("H141").Value = 4180.8335
$ws.Range("I141").Value = 3817
$ws.Range("K141").Value = 11451
$ws.Range("M141").Value = -6271

$ws = $wb.Worksheets.Item("ARM")
# Row 74: Titanium Nugget
$ws.Range("H74").Value = 1423.8334
$ws.Range("I74").Value = 939.4286
$ws.Range("K74").Value = 939.4286
$ws.Range("M74").Value = -65.42859999999996
# Row 77: Titanium Nugget
$ws.Range("H77").Value = 1423.8334
$ws.Range("I77").Value = 939.4286
$ws.Range("K77").Value = 4697.143
$ws.Range("M77").Value = -329.143

$ws = $wb.Worksheets.Item("BSM")
# Row 96: High Steel Sledgehammer
$ws.Range("H96").Value = 19999.5
$ws.Range("I96").Value = 19999.5
$ws.Range("K96").Value = 19999.5
$ws.Range("M96").Value = -17253.5
# Row 99: Oroshigane Ingot
$ws.Range("H99").Value = 1040.6666
$ws.Range("I99").Value = 766.7143
$ws.Range("J99").Value = 1999.5
$ws.Range("K99").Value = 766.7143
$ws.Range("L99").Value = 1999.5
$ws.Range("M99").Value = 731.2857
$ws.Range("N99").Value = -4995.5
# Row 105: Molybdenum Ingot
$ws.Range("H105").Value = 2760.3225
$ws.Range("I105").Value = 2090.92
$ws.Range("J105").Value = 5549.5
$ws.Range("K105").Value = 2090.92
$ws.Range("L105").Value = 5549.5
$ws.Range("M105").Value = -343.9200000000001
$ws.Range("N105").Value = -9043.5

$ws = $wb.Worksheets.Item("CRP")
# Row 7: Maple Lumber
$ws.Range("H7").Value = 118.1
$ws.Range("I7").Value = 98.5
$ws.Range("K7").Value = 98.5
$ws.Range("M7").Value = 14.5
# Row 22: Elm Lumber
$ws.Range("H22").Value = 83703.06
$ws.Range("I22").Value = 119574.8
$ws.Range("K22").Value = 119574.8
$ws.Range("M22").Value = -119224.8
# Row 31: Walnut Lumber
$ws.Range("H31").Value = 3862.6365
$ws.Range("I31").Value = 2967.0667
$ws.Range("K31").Value = 2967.0667
$ws.Range("M31").Value = -2672.0667
# Row 34: Walnut Lumber
$ws.Range("H34").Value = 3862.6365
$ws.Range("I34").Value = 2967.0667
$ws.Range("K34").Value = 2967.0667
$ws.Range("M34").Value = -2765.0667
# Row 109: White Oak Necklace
$ws.Range("H109").Value = 60000
$ws.Range("J109").Value = 60000
$ws.Range("L109").Value = 60000
$ws.Range("N109").Value = -62080
# Row 135: Ceiba Wings
$ws.Range("H135").Value = 125000
$ws.Range("J135").Value = 125000
$ws.Range("L135").Value = 125000
$ws.Range("N135").Value = -135140

$ws = $wb.Worksheets.Item("CUL")
# Row 2: Table Salt
$ws.Range("H2").Value = 342.16666
$ws.Range("I2").Value = 53.75
$ws.Range("J2").Value = 486.375
$ws.Range("K2").Value = 322.5
$ws.Range("L2").Value = 2918.25
$ws.Range("M2").Value = -209.5
$ws.Range("N2").Value = -3144.25
# Row 4: Boiled Egg
$ws.Range("H4").Value = 3812957.2
$ws.Range("J4").Value = 754.5
$ws.Range("L4").Value = 2263.5
$ws.Range("N4").Value = -2487.5
# Row 14: Kukuru Powder
$ws.Range("H14").Value = 573.8889
$ws.Range("I14").Value = 573.8889
$ws.Range("K14").Value = 1721.6667
$ws.Range("M14").Value = -1548.6667
# Row 23: Lavender Oil
$ws.Range("H23").Value = 336667
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 336667
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 1010001
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -1010471
# Row 38: Dark Vinegar
$ws.Range("H38").Value = 106
$ws.Range("J38").Value = 139.33333
$ws.Range("L38").Value = 417.99999
$ws.Range("N38").Value = -1111.99999
# Row 98: Rice Vinegar
$ws.Range("H98").Value = 2187
$ws.Range("I98").Value = 2259.2
$ws.Range("J98").Value = 2066.6667
$ws.Range("K98").Value = 6777.599999999999
$ws.Range("L98").Value = 6200.000100000001
$ws.Range("M98").Value = -5279.599999999999
$ws.Range("N98").Value = -9196.000100000001

$ws = $wb.Worksheets.Item("GSM")
# Row 2: Copper Ingot
$ws.Range("H2").Value = 357.94116
$ws.Range("I2").Value = 95
$ws.Range("K2").Value = 95
$ws.Range("M2").Value = 18
# Row 3: Copper Wristlets
$ws.Range("H3").Value = 3333397
$ws.Range("I3").Value = 3333397
$ws.Range("K3").Value = 3333397
$ws.Range("M3").Value = -3333281
# Row 11: Copper Ring
$ws.Range("H11").Value = 46250000
$ws.Range("I11").Value = 46250000
$ws.Range("K11").Value = 46250000
$ws.Range("M11").Value = -46249861
# Row 18: Brass Gorget
$ws.Range("H18").Value = 510003
$ws.Range("I18").Value = 1000005
$ws.Range("J18").Value = 20001
$ws.Range("K18").Value = 1000005
$ws.Range("L18").Value = 20001
$ws.Range("M18").Value = -999712
$ws.Range("N18").Value = -20587
# Row 80: Hardsilver Ingot
$ws.Range("H80").Value = 4999
$ws.Range("I80").Value = 4999
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 4999
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -4001
$ws.Range("N80").ClearContents()
# Row 83: Hardsilver Ingot
$ws.Range("H83").Value = 4999
$ws.Range("I83").Value = 4999
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 24995
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -20003
$ws.Range("N83").ClearContents()
# Row 126: Phrygian Gold Ingot
$ws.Range("H126").Value = 6085.7144
$ws.Range("J126").Value = 6262.5
$ws.Range("L126").Value = 18787.5
$ws.Range("N126").Value = -23727.5

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Leather
$ws.Range("H7").Value = 1166.3334
$ws.Range("I7").Value = 1166.3334
$ws.Range("K7").Value = 1166.3334
$ws.Range("M7").Value = -1054.3334
# Row 16: Hard Leather
$ws.Range("H16").Value = 1563.3334
$ws.Range("I16").Value = 1350
$ws.Range("K16").Value = 1350
$ws.Range("M16").Value = -1180
# Row 22: Aldgoat Leather
$ws.Range("H22").Value = 2665.5454
$ws.Range("I22").Value = 1708.8182
$ws.Range("J22").Value = 3622.2727
$ws.Range("K22").Value = 1708.8182
$ws.Range("L22").Value = 3622.2727
$ws.Range("M22").Value = -1413.8182
$ws.Range("N22").Value = -4212.2727
# Row 27: Aldgoat Leather
$ws.Range("H27").Value = 2665.5454
$ws.Range("I27").Value = 1708.8182
$ws.Range("J27").Value = 3622.2727
$ws.Range("K27").Value = 1708.8182
$ws.Range("L27").Value = 3622.2727
$ws.Range("M27").Value = -1601.8182
$ws.Range("N27").Value = -3836.2727
# Row 61: Raptor Leather
$ws.Range("H61").Value = 2000
$ws.Range("I61").Value = 2000
$ws.Range("K61").Value = 2000
$ws.Range("M61").Value = -1798
# Row 113: Atrociraptor Leather
$ws.Range("H113").Value = 2000
$ws.Range("I113").Value = 2000
$ws.Range("K113").Value = 2000
$ws.Range("M113").Value = 170
# Row 122: Gaja Leather
$ws.Range("H122").Value = 5671.5
$ws.Range("I122").Value = 3355.125
$ws.Range("J122").Value = 6995.143
$ws.Range("K122").Value = 10065.375
$ws.Range("L122").Value = 20985.429
$ws.Range("M122").Value = -7615.375
$ws.Range("N122").Value = -25885.429
# Row 126: Saiga Leather
$ws.Range("H126").Value = 1166.3334
$ws.Range("I126").Value = 1166.3334
$ws.Range("K126").Value = 3499.0002
$ws.Range("M126").Value = -1029.0002
# Row 127: Saigaskin Coat of Fending
$ws.Range("H127").Value = 42499.5
$ws.Range("J127").Value = 42499.5
$ws.Range("L127").Value = 42499.5
$ws.Range("N127").Value = -52419.5

$ws = $wb.Worksheets.Item("WVR")
# Row 45: Linen Trousers
$ws.Range("H45").Value = 17623.5
$ws.Range("J45").Value = 17623.5
$ws.Range("L45").Value = 17623.5
$ws.Range("N45").Value = -18605.5
